$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1717705
$ws.Range("H2").Value = 0.343541
$ws.Range("I2").Value = 0.04063177891664595
$ws.Range("J2").Value = 0.02745976565347561
$ws.Range("M2").Value = 5.187993000000001
$ws.Range("N2").Value = 10.375986
$ws.Range("O2").Value = 0.2055237229286108
$ws.Range("P2").Value = 0.1647366551777407
$ws.Range("Q2").Value = 0.8911441516065001
$ws.Range("R2").Value = 3.564576606426
$ws.Range("S2").Value = 0.008350794472161315
$ws.Range("T2").Value = 0.00452362994571818
$ws.Range("G3").Value = 0.1717705
$ws.Range("H3").Value = 0.343541
$ws.Range("I3").Value = 0.04063177891664595
$ws.Range("J3").Value = 0.02745976565347561
$ws.Range("M3").Value = 7.026005
$ws.Range("O3").Value = 0.2783370573003924
$ws.Range("P3").Value = 0.3346498047401226
$ws.Range("Q3").Value = 1.2068603918525
$ws.Range("R3").Value = 7.241162351115
$ws.Range("S3").Value = 0.01130932977653936
$ws.Range("T3").Value = 0.00918940521414514
$ws.Range("G4").Value = 0.1717705
$ws.Range("H4").Value = 0.343541
$ws.Range("I4").Value = 0.04063177891664595
$ws.Range("J4").Value = 0.02745976565347561
$ws.Range("M4").Value = 2.488866666666667
$ws.Range("N4").Value = 7.4666
$ws.Range("O4").Value = 0.09859711514766022
$ws.Range("P4").Value = 0.118545139666738
$ws.Range("Q4").Value = 0.4275138717666666
$ws.Range("R4").Value = 2.5650832306
$ws.Range("S4").Value = 0.004006176184498814
$ws.Range("T4").Value = 0.003255221754607163
$ws.Range("G5").Value = 0.1717705
$ws.Range("H5").Value = 0.343541
$ws.Range("I5").Value = 0.04063177891664595
$ws.Range("J5").Value = 0.02745976565347561
$ws.Range("M5").Value = 7.555099
$ws.Range("N5").Value = 15.110198
$ws.Range("O5").Value = 0.299297256872595
$ws.Range("P5").Value = 0.2399004275442726
$ws.Range("Q5").Value = 1.2977431327795
$ws.Range("R5").Value = 5.190972531118
$ws.Range("S5").Value = 0.01216097997160587
$ws.Range("T5").Value = 0.006587609520534333
$ws.Range("G6").Value = 0.1717705
$ws.Range("H6").Value = 0.343541
$ws.Range("I6").Value = 0.04063177891664595
$ws.Range("J6").Value = 0.02745976565347561
$ws.Range("M6").Value = 0.8052820000000001
$ws.Range("N6").Value = 2.415846
$ws.Range("O6").Value = 0.03190146067032041
$ws.Range("P6").Value = 0.03835571766042516
$ws.Range("Q6").Value = 0.138323691781
$ws.Range("R6").Value = 0.829942150686
$ws.Range("S6").Value = 0.001296213097074535
$ws.Range("T6").Value = 0.001053239018426151
$ws.Range("G7").Value = 0.1717705
$ws.Range("H7").Value = 0.343541
$ws.Range("I7").Value = 0.04063177891664595
$ws.Range("J7").Value = 0.02745976565347561
$ws.Range("M7").Value = 2.179548333333333
$ws.Range("N7").Value = 6.538645000000001
$ws.Range("O7").Value = 0.08634338708042119
$ws.Range("P7").Value = 0.1038122552107008
$ws.Range("Q7").Value = 0.3743821069908334
$ws.Range("R7").Value = 2.246292641945
$ws.Range("S7").Value = 0.003508285414766059
$ws.Range("T7").Value = 0.002850660200044646
$ws.Range("I8").Value = 0.9593682210833541
$ws.Range("J8").Value = 0.9725402343465244
$ws.Range("M8").Value = 5.187993000000001
$ws.Range("N8").Value = 10.375986
$ws.Range("O8").Value = 0.2055237229286108
$ws.Range("P8").Value = 0.1647366551777407
$ws.Range("Q8").Value = 21.041052157953
$ws.Range("R8").Value = 126.246312947718
$ws.Range("S8").Value = 0.1971729284564495
$ws.Range("T8").Value = 0.1602130252320225
$ws.Range("I9").Value = 0.9593682210833541
$ws.Range("J9").Value = 0.9725402343465244
$ws.Range("M9").Value = 7.026005
$ws.Range("O9").Value = 0.2783370573003924
$ws.Range("P9").Value = 0.3346498047401226
$ws.Range("S9").Value = 0.267027727523853
$ws.Range("T9").Value = 0.3254603995259775
$ws.Range("I10").Value = 0.9593682210833541
$ws.Range("J10").Value = 0.9725402343465244
$ws.Range("M10").Value = 2.488866666666667
$ws.Range("N10").Value = 7.4666
$ws.Range("O10").Value = 0.09859711514766022
$ws.Range("P10").Value = 0.118545139666738
$ws.Range("Q10").Value = 10.0941488062
$ws.Range("R10").Value = 90.84733925579999
$ws.Range("S10").Value = 0.09459093896316141
$ws.Range("T10").Value = 0.1152899179121309
$ws.Range("I11").Value = 0.9593682210833541
$ws.Range("J11").Value = 0.9725402343465244
$ws.Range("M11").Value = 7.555099
$ws.Range("N11").Value = 15.110198
$ws.Range("O11").Value = 0.299297256872595
$ws.Range("P11").Value = 0.2399004275442726
$ws.Range("Q11").Value = 30.641373671379
$ws.Range("R11").Value = 183.848242028274
$ws.Range("S11").Value = 0.2871362769009891
$ws.Range("T11").Value = 0.2333128180237383
$ws.Range("I12").Value = 0.9593682210833541
$ws.Range("J12").Value = 0.9725402343465244
$ws.Range("M12").Value = 0.8052820000000001
$ws.Range("N12").Value = 2.415846
$ws.Range("O12").Value = 0.03190146067032041
$ws.Range("P12").Value = 0.03835571766042516
$ws.Range("Q12").Value = 3.265999118322001
$ws.Range("R12").Value = 29.393992064898
$ws.Range("S12").Value = 0.03060524757324588
$ws.Range("T12").Value = 0.037302478641999
$ws.Range("I13").Value = 0.9593682210833541
$ws.Range("J13").Value = 0.9725402343465244
$ws.Range("M13").Value = 2.179548333333333
$ws.Range("N13").Value = 6.538645000000001
$ws.Range("O13").Value = 0.08634338708042119
$ws.Range("P13").Value = 0.3346498047401226
$ws.Range("Q13").Value = 8.839639946015001
$ws.Range("R13").Value = 79.55675951413501
$ws.Range("S13").Value = 0.08283510166565514
$ws.Range("T13").Value = 0.1009615950106561
